$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.73885133333334
$ws.Range("H2").Value = 107.216554
$ws.Range("I2").Value = 0.01949729408921566
$ws.Range("J2").Value = 0.01949729408921566
$ws.Range("M2").Value = 24.576554
$ws.Range("N2").Value = 73.729662
$ws.Range("O2").Value = 0.07553767049546639
$ws.Range("P2").Value = 0.07553767049546638
$ws.Range("Q2").Value = 878.3378096916388
$ws.Range("R2").Value = 7905.040287224749
$ws.Range("S2").Value = 0.001472780176464377
$ws.Range("T2").Value = 0.001472780176464377
$ws.Range("G3").Value = 35.73885133333334
$ws.Range("H3").Value = 107.216554
$ws.Range("I3").Value = 0.01949729408921566
$ws.Range("J3").Value = 0.01949729408921566
$ws.Range("O3").Value = 0.359764849016532
$ws.Range("P3").Value = 0.359764849016532
$ws.Range("Q3").Value = 4183.27792499491
$ws.Range("R3").Value = 37649.50132495419
$ws.Range("S3").Value = 0.007014441064237594
$ws.Range("T3").Value = 0.007014441064237592
$ws.Range("G4").Value = 35.73885133333334
$ws.Range("H4").Value = 107.216554
$ws.Range("I4").Value = 0.01949729408921566
$ws.Range("J4").Value = 0.01949729408921566
$ws.Range("M4").Value = 55.68784966666667
$ws.Range("N4").Value = 167.063549
$ws.Range("O4").Value = 0.1711603033819035
$ws.Range("P4").Value = 0.1711603033819035
$ws.Range("Q4").Value = 1990.219780310017
$ws.Range("R4").Value = 17911.97802279015
$ws.Range("S4").Value = 0.003337162771436347
$ws.Range("T4").Value = 0.003337162771436346
$ws.Range("G5").Value = 35.73885133333334
$ws.Range("H5").Value = 107.216554
$ws.Range("I5").Value = 0.01949729408921566
$ws.Range("J5").Value = 0.01949729408921566
$ws.Range("M5").Value = 128.0392633333333
$ws.Range("N5").Value = 384.11779
$ws.Range("O5").Value = 0.3935371771060981
$ws.Range("P5").Value = 0.3935371771060981
$ws.Range("Q5").Value = 4575.976197099519
$ws.Range("R5").Value = 41183.78577389566
$ws.Range("S5").Value = 0.007672910077077343
$ws.Range("T5").Value = 0.007672910077077342
$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.9215900675332435
$ws.Range("J6").Value = 0.9215900675332435
$ws.Range("M6").Value = 24.576554
$ws.Range("N6").Value = 73.729662
$ws.Range("O6").Value = 0.07553767049546639
$ws.Range("P6").Value = 0.07553767049546638
$ws.Range("Q6").Value = 41516.90986691589
$ws.Range("R6").Value = 373652.188802243
$ws.Range("S6").Value = 0.06961476685322077
$ws.Range("T6").Value = 0.06961476685322075
$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.9215900675332435
$ws.Range("J7").Value = 0.9215900675332435
$ws.Range("O7").Value = 0.359764849016532
$ws.Range("P7").Value = 0.359764849016532
$ws.Range("S7").Value = 0.3315557115012329
$ws.Range("T7").Value = 0.3315557115012329
$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.9215900675332435
$ws.Range("J8").Value = 0.9215900675332435
$ws.Range("M8").Value = 55.68784966666667
$ws.Range("N8").Value = 167.063549
$ws.Range("O8").Value = 0.1711603033819035
$ws.Range("P8").Value = 0.1711603033819035
$ws.Range("Q8").Value = 94072.8889531609
$ws.Range("R8").Value = 846656.0005784482
$ws.Range("S8").Value = 0.1577396355527389
$ws.Range("T8").Value = 0.1577396355527389
$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.9215900675332435
$ws.Range("J9").Value = 0.9215900675332435
$ws.Range("M9").Value = 128.0392633333333
$ws.Range("N9").Value = 384.11779
$ws.Range("O9").Value = 0.3935371771060981
$ws.Range("P9").Value = 0.3935371771060981
$ws.Range("Q9").Value = 216295.3583824774
$ws.Range("R9").Value = 1946658.225442297
$ws.Range("S9").Value = 0.3626799536260509
$ws.Range("T9").Value = 0.3626799536260509
$ws.Range("G10").Value = 93.641553
$ws.Range("H10").Value = 280.924659
$ws.Range("I10").Value = 0.05108605424341119
$ws.Range("J10").Value = 0.05108605424341119
$ws.Range("M10").Value = 24.576554
$ws.Range("N10").Value = 73.729662
$ws.Range("O10").Value = 0.07553767049546639
$ws.Range("P10").Value = 0.07553767049546638
$ws.Range("Q10").Value = 2301.386683948362
$ws.Range("R10").Value = 20712.48015553526
$ws.Range("S10").Value = 0.003858921532352317
$ws.Range("T10").Value = 0.003858921532352316
$ws.Range("G11").Value = 93.641553
$ws.Range("H11").Value = 280.924659
$ws.Range("I11").Value = 0.05108605424341119
$ws.Range("J11").Value = 0.05108605424341119
$ws.Range("O11").Value = 0.359764849016532
$ws.Range("P11").Value = 0.359764849016532
$ws.Range("Q11").Value = 10960.86267220846
$ws.Range("R11").Value = 98647.76404987615
$ws.Range("S11").Value = 0.01837896659173119
$ws.Range("T11").Value = 0.01837896659173119
$ws.Range("G12").Value = 93.641553
$ws.Range("H12").Value = 280.924659
$ws.Range("I12").Value = 0.05108605424341119
$ws.Range("J12").Value = 0.05108605424341119
$ws.Range("M12").Value = 55.68784966666667
$ws.Range("N12").Value = 167.063549
$ws.Range("O12").Value = 0.1711603033819035
$ws.Range("P12").Value = 0.1711603033819035
$ws.Range("Q12").Value = 5214.6967260172
$ws.Range("R12").Value = 46932.2705341548
$ws.Range("S12").Value = 0.008743904542886639
$ws.Range("T12").Value = 0.008743904542886639
$ws.Range("G13").Value = 93.641553
$ws.Range("H13").Value = 280.924659
$ws.Range("I13").Value = 0.05108605424341119
$ws.Range("J13").Value = 0.05108605424341119
$ws.Range("M13").Value = 128.0392633333333
$ws.Range("N13").Value = 384.11779
$ws.Range("O13").Value = 0.3935371771060981
$ws.Range("P13").Value = 0.3935371771060981
$ws.Range("Q13").Value = 11989.79546350929
$ws.Range("R13").Value = 107908.1591715836
$ws.Range("S13").Value = 0.02010426157644104
$ws.Range("T13").Value = 0.02010426157644104
$ws.Range("G14").Value = 14.34625366666667
$ws.Range("H14").Value = 43.038761
$ws.Range("I14").Value = 0.007826584134129748
$ws.Range("J14").Value = 0.007826584134129748
$ws.Range("M14").Value = 24.576554
$ws.Range("N14").Value = 73.729662
$ws.Range("O14").Value = 0.07553767049546639
$ws.Range("P14").Value = 0.07553767049546638
$ws.Range("Q14").Value = 352.5814779365314
$ws.Range("R14").Value = 3173.233301428782
$ws.Range("S14").Value = 0.0005912019334289381
$ws.Range("T14").Value = 0.000591201933428938
$ws.Range("G15").Value = 14.34625366666667
$ws.Range("H15").Value = 43.038761
$ws.Range("I15").Value = 0.007826584134129748
$ws.Range("J15").Value = 0.007826584134129748
$ws.Range("O15").Value = 0.359764849016532
$ws.Range("P15").Value = 0.359764849016532
$ws.Range("Q15").Value = 1679.247206643406
$ws.Range("R15").Value = 15113.22485979065
$ws.Range("S15").Value = 0.002815729859330374
$ws.Range("T15").Value = 0.002815729859330374
$ws.Range("G16").Value = 14.34625366666667
$ws.Range("H16").Value = 43.038761
$ws.Range("I16").Value = 0.007826584134129748
$ws.Range("J16").Value = 0.007826584134129748
$ws.Range("M16").Value = 55.68784966666667
$ws.Range("N16").Value = 167.063549
$ws.Range("O16").Value = 0.1711603033819035
$ws.Range("P16").Value = 0.1711603033819035
$ws.Range("Q16").Value = 798.9120174691989
$ws.Range("R16").Value = 7190.20815722279
$ws.Range("S16").Value = 0.00133960051484164
$ws.Range("T16").Value = 0.00133960051484164
$ws.Range("G17").Value = 14.34625366666667
$ws.Range("H17").Value = 43.038761
$ws.Range("I17").Value = 0.007826584134129748
$ws.Range("J17").Value = 0.007826584134129748
$ws.Range("M17").Value = 128.0392633333333
$ws.Range("N17").Value = 384.11779
$ws.Range("O17").Value = 0.3935371771060981
$ws.Range("P17").Value = 0.3935371771060981
$ws.Range("Q17").Value = 1836.883751073132
$ws.Range("R17").Value = 16531.95375965819
$ws.Range("S17").Value = 0.003080051826528796
$ws.Range("T17").Value = 0.003080051826528796
